# Rename the sheet formerly known as "Property1" to "DataNode" as part of
# unifying the DataNode / DataTable / Entity naming convention across the
# data-config workbooks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "DataNode"

# Reflect the last on-screen selection that was active when the workbook
# was saved (bottom-left frozen pane moved from K9 to O40).
[void]$ws.Range("O40").Select()
